$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) columns
# Numeric-looking text values need NumberFormat "@" forced so Excel
# keeps them as text (preserving trailing zeros / multi-dot formatting),
# then the cell style is reset to Normal so no stray format lingers.

$ws.Range("D2").Value = '58.920.39'
$ws.Range("E2").Value = '  -0.82%  '
$ws.Range("D3").Value = '2.606.42'
$ws.Range("E3").Value = '  -1.24%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '557.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.99%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("E8").Value = '  +4.19%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.83'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.93%  '
$ws.Range("E10").Value = '  -0.76%  '
$ws.Range("E11").Value = '  +5.91%  '
$ws.Range("E12").Value = '  -0.92%  '
$ws.Range("D13").Value = '3.065.20'
$ws.Range("E13").Value = '  -1.33%  '
$ws.Range("D14").Value = '58.875.40'
$ws.Range("E14").Value = '  -0.77%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.02'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.22%  '
$ws.Range("D16").Value = '2.618.29'
$ws.Range("E16").Value = '  -0.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000133'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.71%  '
$ws.Range("E18").Value = '  -0.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '338.01'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.12'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.73%  '
$ws.Range("E21").Value = '  -0.80%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.25'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.429'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.96%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.22%  '
$ws.Range("E26").Value = '  -2.26%  '
$ws.Range("E27").Value = '  -2.13%  '
$ws.Range("D28").Value = '0.0₃0762'
$ws.Range("E28").Value = '  +1.11%  '
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("E30").Value = '  +1.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '154.41'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.37%  '
$ws.Range("E33").Value = '  +0.48%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.97'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.908'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.28%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.906'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.80%  '
$ws.Range("E37").Value = '  -0.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '37.16'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.54%  '
$ws.Range("E39").Value = '  +1.34%  '
$ws.Range("E40").Value = '  -0.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '283.42'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.72%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.997'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.600'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.10%  '
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0955'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.14%  '
$ws.Range("E46").Value = '  -1.34%  '
$ws.Range("E47").Value = '  +2.55%  '
$ws.Range("E48").Value = '  +0.06%  '
$ws.Range("D49").Value = '1.943.66'
$ws.Range("E49").Value = '  -0.96%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '117.76'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.98'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.36%  '
